$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()

$ws.Range("H19").Value = 953.5714
$ws.Range("I19").Value = 980.8333
$ws.Range("K19").Value = 980.8333
$ws.Range("M19").Value = -805.8333

$ws.Range("H31").Value = 39.8
$ws.Range("I31").Value = 39.8
$ws.Range("K31").Value = 119.4
$ws.Range("M31").Value = 110.6

$ws.Range("H70").Value = 3758.84
$ws.Range("J70").Value = 4155.8096
$ws.Range("L70").Value = 12467.4288
$ws.Range("N70").Value = -13007.4288

$ws.Range("H73").Value = 3758.84
$ws.Range("J73").Value = 4155.8096
$ws.Range("L73").Value = 12467.4288
$ws.Range("N73").Value = -14339.4288

$ws.Range("H98").Value = 1053.5
$ws.Range("I98").Value = 1225.1333
$ws.Range("K98").Value = 1225.1333
$ws.Range("M98").Value = 272.8667

$ws.Range("H116").Value = 4400
$ws.Range("I116").Value = 3800
$ws.Range("K116").Value = 3800
$ws.Range("M116").Value = -358

$ws.Range("H122").Value = 1053.5
$ws.Range("I122").Value = 1225.1333
$ws.Range("K122").Value = 3675.3999
$ws.Range("M122").Value = -1225.3999

$ws.Range("H132").Value = 3058.8667
$ws.Range("I132").Value = 3058.8667
$ws.Range("K132").Value = 9176.6001
$ws.Range("M132").Value = -6646.6001

$ws.Range("H137").Value = 2016.0571
$ws.Range("I137").Value = 1594.2693
$ws.Range("J137").Value = 3234.5557
$ws.Range("K137").Value = 4782.8079
$ws.Range("L137").Value = 9703.667099999999
$ws.Range("M137").Value = -2232.8079
$ws.Range("N137").Value = -14803.6671

$ws.Range("H138").Value = 3555.0232
$ws.Range("J138").Value = 2660.7778
$ws.Range("L138").Value = 7982.3334
$ws.Range("N138").Value = -18262.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 5999
$ws.Range("I31").Value = 5999
$ws.Range("K31").Value = 5999
$ws.Range("M31").Value = -5705

$ws.Range("H32").Value = 5732.174
$ws.Range("I32").Value = 3820.3508
$ws.Range("J32").Value = 14813.333
$ws.Range("K32").Value = 3820.3508
$ws.Range("L32").Value = 14813.333
$ws.Range("M32").Value = -3533.3508
$ws.Range("N32").Value = -15387.333

$ws.Range("H61").Value = 1514.12
$ws.Range("I61").Value = 1411.0869
$ws.Range("K61").Value = 1411.0869
$ws.Range("M61").Value = -1199.0869

$ws.Range("H102").Value = 1553.9
$ws.Range("I102").Value = 1442
$ws.Range("K102").Value = 1442
$ws.Range("M102").Value = 180

$ws.Range("H132").Value = 2859.4443
$ws.Range("I132").Value = 2779.375
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 8338.125
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -5808.125
$ws.Range("N132").Value = -15560

$ws.Range("H136").Value = 1514.12
$ws.Range("I136").Value = 1411.0869
$ws.Range("K136").Value = 4233.2607
$ws.Range("M136").Value = -1683.2607

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

$ws.Range("H80").Value = 879.8333
$ws.Range("I80").Value = 912
$ws.Range("J80").Value = 719
$ws.Range("K80").Value = 912
$ws.Range("L80").Value = 719
$ws.Range("M80").Value = 86
$ws.Range("N80").Value = -2715

$ws.Range("H83").Value = 879.8333
$ws.Range("I83").Value = 912
$ws.Range("J83").Value = 719
$ws.Range("K83").Value = 4560
$ws.Range("L83").Value = 3595
$ws.Range("M83").Value = 432
$ws.Range("N83").Value = -13579

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3236.64
$ws.Range("I58").Value = 1162.909
$ws.Range("K58").Value = 1162.909
$ws.Range("M58").Value = -959.9090000000001

$ws.Range("H103").Value = 14898.333
$ws.Range("I103").Value = 14898.333
$ws.Range("K103").Value = 14898.333
$ws.Range("M103").Value = -13726.333

$ws.Range("H105").Value = 1170.25
$ws.Range("I105").Value = 890
$ws.Range("K105").Value = 890
$ws.Range("M105").Value = 857

$ws.Range("H132").Value = 1973.1562
$ws.Range("J132").Value = 1185.2
$ws.Range("L132").Value = 3555.6
$ws.Range("N132").Value = -8615.6

$ws.Range("H134").Value = 1819.5385
$ws.Range("I134").Value = 876.3333
$ws.Range("J134").Value = 3941.75
$ws.Range("K134").Value = 2628.9999
$ws.Range("L134").Value = 11825.25
$ws.Range("M134").Value = -93.9998999999998
$ws.Range("N134").Value = -16895.25

$ws.Range("H136").Value = 3236.64
$ws.Range("I136").Value = 1162.909
$ws.Range("K136").Value = 3488.727
$ws.Range("M136").Value = -938.7270000000003

$ws.Range("H141").Value = 133853.42
$ws.Range("J141").Value = 149875.6
$ws.Range("L141").Value = 149875.6
$ws.Range("N141").Value = -160235.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 83334060
$ws.Range("I64").Value = 112
$ws.Range("J64").Value = 111112050
$ws.Range("K64").Value = 336
$ws.Range("L64").Value = 333336150
$ws.Range("M64").Value = -66
$ws.Range("N64").Value = -333336690

$ws.Range("H67").Value = 83334060
$ws.Range("I67").Value = 112
$ws.Range("J67").Value = 111112050
$ws.Range("K67").Value = 336
$ws.Range("L67").Value = 333336150
$ws.Range("M67").Value = 600
$ws.Range("N67").Value = -333338022

$ws.Range("H110").Value = 598.5
$ws.Range("I110").Value = 598.5
$ws.Range("K110").Value = 1795.5
$ws.Range("M110").Value = 2294.5

$ws.Range("H111").Value = 485.75
$ws.Range("I111").Value = 485.75
$ws.Range("K111").Value = 1457.25
$ws.Range("M111").Value = 1609.75

$ws.Range("H118").Value = 936.0769
$ws.Range("I118").Value = 681.2857
$ws.Range("J118").Value = 1233.3334
$ws.Range("K118").Value = 2043.8571
$ws.Range("L118").Value = 3700.0002
$ws.Range("M118").Value = -800.8571000000002
$ws.Range("N118").Value = -6186.0002

$ws.Range("H120").Value = 14730.77
$ws.Range("I120").Value = 10500
$ws.Range("J120").Value = 16000
$ws.Range("K120").Value = 31500
$ws.Range("L120").Value = 48000
$ws.Range("M120").Value = -26662
$ws.Range("N120").Value = -57676

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 87000
$ws.Range("J103").Value = 87000
$ws.Range("L103").Value = 87000
$ws.Range("N103").Value = -89344

$ws.Range("H132").Value = 5230.8335
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3013.4666
$ws.Range("I46").Value = 2666.6667
$ws.Range("J46").Value = 3244.6667
$ws.Range("K46").Value = 2666.6667
$ws.Range("L46").Value = 3244.6667
$ws.Range("M46").Value = -2478.6667
$ws.Range("N46").Value = -3620.6667

$ws.Range("H100").Value = 2033.3334
$ws.Range("I100").Value = 1100
$ws.Range("K100").Value = 1100
$ws.Range("M100").Value = -559

$ws.Range("H122").Value = 4383.625
$ws.Range("I122").Value = 3300.1667
$ws.Range("J122").Value = 5033.7
$ws.Range("K122").Value = 9900.500100000001
$ws.Range("L122").Value = 15101.1
$ws.Range("M122").Value = -7450.500100000001
$ws.Range("N122").Value = -20001.1

$ws.Range("H136").Value = 1780.25
$ws.Range("I136").Value = 1714.5526
$ws.Range("K136").Value = 5143.6578
$ws.Range("M136").Value = -2593.6578

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

$ws.Range("H132").Value = 1110.5625
$ws.Range("I132").Value = 1155.4667
$ws.Range("K132").Value = 3466.4001
$ws.Range("M132").Value = -936.4000999999998
